$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94 (shifts existing rows 94-97 down to 95-98)
$ws.Rows.Item(94).Insert()

# Populate the new row 94 with this week's data
$ws.Cells.Item(94, 1).Value = 6
$ws.Cells.Item(94, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(94, 3).Value = "Metropolitana"
$ws.Cells.Item(94, 4).Value = 44461
$ws.Cells.Item(94, 5).Value = 13
$ws.Cells.Item(94, 6).Value = 100112001
$ws.Cells.Item(94, 7).Value = "Berenjena"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 580
$ws.Cells.Item(94, 11).Value = 6000
$ws.Cells.Item(94, 12).Value = 7000
$ws.Cells.Item(94, 13).Value = 6552
$ws.Cells.Item(94, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(94, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(94, 16).Value = 131
$ws.Cells.Item(94, 17).Value = 50
$ws.Cells.Item(94, 18).Value = "Hortaliza"

# Match the date-style formatting used by column D (style index 2 / custom date format)
$ws.Cells.Item(94, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
